$wb = $excel.ActiveWorkbook

# Sheet 1: Exhibition 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 238
$ws.Range("F4").Value = 496
$ws.Range("F5").Value = 2162
$ws.Range("F7").Value = 8099
$ws.Range("F8").Value = 287
$ws.Range("F9").Value = 75
$ws.Range("F11").Value = 255
$ws.Range("F12").Value = 1778
$ws.Range("F13").Value = 1559
$ws.Range("F14").Value = 1317
$ws.Range("F15").Value = 194
$ws.Range("F16").Value = 4169
$ws.Range("F17").Value = 715
$ws.Range("F18").Value = 43
$ws.Range("F19").Value = 1145
$ws.Range("F20").Value = 1240
$ws.Range("F21").Value = 441
$ws.Range("F22").Value = 6284
$ws.Range("F25").Value = 4265
$ws.Range("F27").Value = 1964
$ws.Range("F28").Value = 1173
$ws.Range("F30").Value = 1045
$ws.Range("F31").Value = 25
$ws.Range("F33").Value = 212
$ws.Range("F34").Value = 58
$ws.Range("F35").Value = 327
$ws.Range("F36").Value = 1164
$ws.Range("F37").Value = 515
$ws.Range("F38").Value = 1885
$ws.Range("F39").Value = 123
$ws.Range("F40").Value = 422
$ws.Range("F42").Value = 1167
$ws.Range("F43").Value = 560
$ws.Range("F45").Value = 961

# Sheet 2: Performance 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 992
$ws.Range("F5").Value = 123
$ws.Range("F9").Value = 2
$ws.Range("F11").Value = 18
$ws.Range("F14").Value = 383
$ws.Range("F15").Value = 407
$ws.Range("F17").Value = 206
$ws.Range("F18").Value = 109
$ws.Range("F20").Value = 6
$ws.Range("F24").Value = 97
$ws.Range("F27").Value = 96
$ws.Range("F29").Value = 127
$ws.Range("F33").Value = 37
$ws.Range("F35").Value = 6

# Sheet 3: Local Life 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F6").Value = 1574
$ws.Range("F7").Value = 477
$ws.Range("F8").Value = 3117
$ws.Range("F9").Value = 972
$ws.Range("F10").Value = 1113
$ws.Range("F11").Value = 1315
$ws.Range("F12").Value = 1612
$ws.Range("F13").Value = 19

# Sheet 4: All Types 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 1574
$ws.Range("F4").Value = 496
$ws.Range("F5").Value = 477
$ws.Range("F6").Value = 3117
$ws.Range("F7").Value = 8099
$ws.Range("F8").Value = 287
$ws.Range("F9").Value = 75
$ws.Range("F10").Value = 972
$ws.Range("F12").Value = 18
$ws.Range("F13").Value = 255
$ws.Range("F14").Value = 1778
$ws.Range("F15").Value = 1560
$ws.Range("F16").Value = 1317
$ws.Range("F18").Value = 194
$ws.Range("F19").Value = 1612
$ws.Range("F20").Value = 4169
$ws.Range("F21").Value = 383
$ws.Range("F22").Value = 407
$ws.Range("F23").Value = 715
$ws.Range("F24").Value = 43
$ws.Range("F25").Value = 1145
$ws.Range("F26").Value = 1240
$ws.Range("F27").Value = 441
$ws.Range("F28").Value = 6284
$ws.Range("F31").Value = 1964
$ws.Range("F32").Value = 1173
$ws.Range("F34").Value = 25
$ws.Range("F35").Value = 212
$ws.Range("F36").Value = 58
$ws.Range("F37").Value = 97
$ws.Range("F38").Value = 327
$ws.Range("F39").Value = 1164
$ws.Range("F40").Value = 515
$ws.Range("F41").Value = 1885
$ws.Range("F42").Value = 123
$ws.Range("F43").Value = 422
$ws.Range("F44").Value = 1167
$ws.Range("F45").Value = 127
$ws.Range("F46").Value = 560
